$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44308, 4, 24, 238.3316782522344),
    @(44309, 2, 23, 228.4011916583913),
    @(44310, 4, 23, 228.4011916583913),
    @(44311, 5, 26, 258.1926514399206),
    @(44312, 7, 30, 297.914597815293)
)

$startRow = 234

# Copy the formatting of the last existing data row (233) down onto the
# new rows first, so the new cells inherit the same style (esp. column A's
# date style) without introducing new cellXfs entries.
$ws.Range("A233:D233").Copy()
$ws.Range("A234:D238").PasteSpecial(-4122)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
